$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15. This shifts the existing rows 15-26
# down to 16-27, preserving all of their data/formatting.
$ws.Rows.Item(15).Insert()

# Copy the date cell style (s="2") from the row now below (row 16) onto
# the new row's date cell so it keeps the same date number format.
$ws.Range("D16").Copy()
$ws.Range("D15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row 15 with the new weekly record.
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C15").Value = "Los Lagos"
$ws.Range("D15").Value = (Get-Date -Year 2022 -Month 1 -Day 25).Date
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 100112030
$ws.Range("G15").Value = "Poroto granado"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 31000
$ws.Range("L15").Value = 31000
$ws.Range("M15").Value = 31000
$ws.Range("N15").Value = "`$/saco 25 kilos"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 1240
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
